$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login_test")

# Update header labels to lowercase wording (locators)
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "password"

# Keep data row values as-is (re-assert to be safe)
$ws.Range("A2").Value = "admin@yourstore.com"
$ws.Range("B2").Value = "admin"

# Update the selected cell/range shown in the sheet view
$ws.Range("C7").Select()
